$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ST_CHECKOUT_02) text updates ---
$ws.Range("B2").Value = "Bỏ trống thông tin bắt buộc"
$ws.Range("C2").Value = "1. Vào Giỏ -> Thanh toán`n2. Để trống Address & Phone`n3. Click Đặt hàng"
$ws.Range("E2").Value = "Vẫn ở trang checkout do trình duyệt/HTML5 validation chặn submit"

# --- Row 3 (ST_CHECKOUT_01) text updates ---
$ws.Range("B3").Value = "Đặt hàng thành công"
$ws.Range("C3").Value = "1. Vào Giỏ -> Thanh toán`n2. Điền thông tin Address & Phone hợp lệ`n3. Submit & Check Success"
$ws.Range("D3").Value = "Address: HCM, Phone: 0987654321"
$ws.Range("E3").Value = "Chuyển về trang chủ (view-products) và hiển thị thông báo thành công"
$ws.Range("F3").Value = "URL: http://localhost:8080/ShopDuck/user/view-products?success=true"

# --- New row 4 (ST_CHECKOUT_03) ---
$ws.Range("A4").Value = "ST_CHECKOUT_03"
$ws.Range("B4").Value = "Số điện thoại sai định dạng"
$ws.Range("C4").Value = "1. Nhập SĐT là chữ 'abc'`n2. Click Đặt hàng`n3. Check validation"
$ws.Range("D4").Value = "Phone: abc_khong_phai_so"
$ws.Range("E4").Value = "Không cho submit/Vẫn ở trang checkout"
$ws.Range("F4").Value = "URL: http://localhost:8080/ShopDuck/user/checkout"
$ws.Range("G2").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("G4").Value = "PASS"

# --- New row 5 (ST_CHECKOUT_04) ---
$ws.Range("A5").Value = "ST_CHECKOUT_04"
$ws.Range("B5").Value = "Quay lại giỏ hàng"
$ws.Range("C5").Value = "1. Vào trang Thanh toán`n2. Click link/nút Quay lại Giỏ hàng`n3. Check URL"
$ws.Range("D5").Value = "Action: Back to Cart"
$ws.Range("E5").Value = "Về trang view-cart.jsp"
$ws.Range("F5").Value = "URL: http://localhost:8080/ShopDuck/user/order/view-cart.jsp"
$ws.Range("G2").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G5").Value = "PASS"

# --- New row 6 (ST_CHECKOUT_05) ---
$ws.Range("A6").Value = "ST_CHECKOUT_05"
$ws.Range("B6").Value = "Truy cập Checkout khi giỏ rỗng"
$ws.Range("C6").Value = "1. Xóa hết hàng trong giỏ`n2. Truy cập thẳng URL /checkout`n3. Check chuyển hướng"
$ws.Range("D6").Value = "Cart: Empty"
$ws.Range("E6").Value = "Bị đá về trang sản phẩm hoặc giỏ hàng (Không được phép vào trang checkout)"
$ws.Range("F6").Value = "URL: http://localhost:8080/ShopDuck/user/order/view-cart.jsp"
$ws.Range("G2").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G6").Value = "PASS"

# Reset auto-calculated row heights so newly written multi-line rows
# don't pick up an explicit customHeight (match plain rows 2/3 above).
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()

# Re-size (bestFit) columns B-F for the new/expanded content; column widths
# snap to the nearest 1/6 character unit through this API, so these are the
# closest achievable values to the target bestFit widths.
$ws.Columns.Item(2).ColumnWidth = 27.666666666666664
$ws.Columns.Item(3).ColumnWidth = 37.16666666666667
$ws.Columns.Item(4).ColumnWidth = 31.833333333333336
$ws.Columns.Item(5).ColumnWidth = 70.0
$ws.Columns.Item(6).ColumnWidth = 64.66666666666666
